# Update the live cryptocurrency price/volume snapshot on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.012.64'
$ws.Range("E2").Value = '  +9.32%  '
$ws.Range("D3").Value = '2.679.85'
$ws.Range("E3").Value = '  +10.20%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '187.43'
$ws.Range("E5").Value = '  +12.40%  '
$ws.Range("D6").Value = '588.41'
$ws.Range("E6").Value = '  +4.31%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").Value = '0.539'
$ws.Range("E8").Value = '  +4.83%  '
$ws.Range("D9").Value = '0.195'
$ws.Range("E9").Value = '  +15.00%  '
$ws.Range("D10").Value = '2.675.32'
$ws.Range("E10").Value = '  +10.08%  '
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").Value = '0.359'
$ws.Range("E12").Value = '  +6.96%  '
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = '75.670.83'
$ws.Range("E14").Value = '  +9.19%  '
$ws.Range("D15").Value = '3.170.19'
$ws.Range("E15").Value = '  +10.10%  '
$ws.Range("E16").Value = '  +5.91%  '
$ws.Range("D17").Value = '26.61'
$ws.Range("E17").Value = '  +10.95%  '
$ws.Range("D18").Value = '2.680.77'
$ws.Range("E18").Value = '  +11.04%  '
$ws.Range("D19").Value = '9.39'
$ws.Range("E19").Value = '  +31.07%  '
$ws.Range("D20").Value = '11.99'
$ws.Range("E20").Value = '  +10.77%  '
$ws.Range("D21").Value = '373.93'
$ws.Range("E21").Value = '  +9.19%  '
$ws.Range("E22").Value = '  +16.51%  '
$ws.Range("E23").Value = '  +5.05%  '
$ws.Range("E24").Value = '  +4.30%  '
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").Value = '70.18'
$ws.Range("E26").Value = '  +6.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.20'
$ws.Range("E27").Value = '  +9.58%  '
$ws.Range("E28").Value = '  +10.03%  '
$ws.Range("D29").Value = '2.816.26'
$ws.Range("E29").Value = '  +10.23%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").Value = '0.0₃0945'
$ws.Range("E31").Value = '  +11.25%  '
$ws.Range("D32").Value = '523.73'
$ws.Range("E32").Value = '  +15.56%  '
$ws.Range("E33").Value = '  +14.07%  '
$ws.Range("D34").Value = '7.77'
$ws.Range("E34").Value = '  +4.89%  '
$ws.Range("D35").Value = '1.76'
$ws.Range("E35").Value = '  +9.00%  '
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").Value = '163.86'
$ws.Range("E37").Value = '  +3.84%  '
$ws.Range("E38").Value = '  +8.11%  '
$ws.Range("D39").Value = '19.25'
$ws.Range("E39").Value = '  +5.49%  '
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("E42").Value = '  +13.85%  '
$ws.Range("D43").Value = '170.29'
$ws.Range("E43").Value = '  +26.03%  '
$ws.Range("E44").Value = '  +12.32%  '
$ws.Range("D45").Value = '0.331'
$ws.Range("E45").Value = '  +8.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.20'
$ws.Range("E46").Value = '  +9.80%  '
$ws.Range("D47").Value = '39.32'
$ws.Range("D48").Value = '2.36'
$ws.Range("E48").Value = '  +12.72%  '
$ws.Range("D49").Value = '0.0841'
$ws.Range("E49").Value = '  +16.13%  '
$ws.Range("E50").Value = '  +7.29%  '
$ws.Range("D51").Value = '0.538'
$ws.Range("E51").Value = '  +10.07%  '
